# Generate Report for Handoff
#
# Adds a new "localization status" record for file
#   ffff5d5a0c9b-2eae-4d19-adc7-c8a921d9a6c9.md
# as row 3 on every sheet, and refreshes the existing row-2 record for
# 6cbd9040-2d4e-4c40-8a58-d6278cda3a29.md -> d3f57a95-35c5-4236-bd79-77497563d9a8.md
# (new handoff/target filenames + timestamps).

$wb = $excel.ActiveWorkbook

$oldBase = "6cbd9040-2d4e-4c40-8a58-d6278cda3a29"
$newBase = "d3f57a95-35c5-4236-bd79-77497563d9a8"
$newFile2Base = "ffff5d5a0c9b-2eae-4d19-adc7-c8a921d9a6c9"

$newMdName = "$newBase.md"
$newFile2MdName = "$newFile2Base.md"

$newHandoffDate = "2016-03-24 02:55:11"
$zhHandoffDatetime = "2016-03-24 02:55:02"

$zhXlf = "$newBase.05004d211646c709bd999f80776f1e21d62bc427.zh-cn.xlf"
$deXlf = "$newBase.05004d211646c709bd999f80776f1e21d62bc427.de-de.xlf"

$hyperlinkUnderline = 2
$hyperlinkColor = 15570276   # BGR for RGB FF6495ED (matches the workbook's existing HyperLink font colour)

function Set-HyperlinkLook($range) {
    $range.Font.Underline = $hyperlinkUnderline
    $range.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Row 2: point the existing overview record at the new handoff file + date.
$ws1.Range("B2").Value2 = "Ready for handoff"
$ws1.Range("C2").Value2 = "Ready for handoff"
$ws1.Range("D2").Value2 = $newHandoffDate
$ws1.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Row 3: brand-new overview record.
$ws1.Range("A3").Value2 = $newFile2MdName
$ws1.Range("B3").Value2 = "Ready for handoff"
$ws1.Range("C3").Value2 = "Ready for handoff"
$ws1.Range("D3").Value2 = $newHandoffDate
$ws1.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Hyperlinks: wipe the sheet's hyperlinks once (Range.Hyperlinks.Delete clears
# the whole worksheet collection in this engine) then re-add every link fresh
# so relationship ids come out sequential again.
$ws1.Range("A2").Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/bea602d80ca40d9304fa3e04a2adf05645224102/e2e/$newMdName", "", "", $newMdName)
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/bea602d80ca40d9304fa3e04a2adf05645224102/e2e/$newFile2MdName", "", "", $newFile2MdName)
Set-HyperlinkLook $ws1.Range("A2")
Set-HyperlinkLook $ws1.Range("A3")

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("B2").Value2 = ".md"
$ws2.Range("C2").Value2 = "Ready for handoff"
$ws2.Range("D2").Value2 = $zhXlf
$ws2.Range("E2").Value2 = $zhHandoffDatetime
$ws2.Range("H2").Value2 = "0001-01-01 00:00:00"
$ws2.Range("J2").Value2 = "Include"
$ws2.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws2.Range("A3").Value2 = $newFile2MdName
$ws2.Range("B3").Value2 = ".md"
$ws2.Range("C3").Value2 = "Ready for handoff"
$ws2.Range("D3").Value2 = $zhXlf
$ws2.Range("E3").Value2 = $zhHandoffDatetime
$ws2.Range("H3").Value2 = "0001-01-01 00:00:00"
$ws2.Range("J3").Value2 = "Include"
$ws2.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws2.Range("A2").Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/bea602d80ca40d9304fa3e04a2adf05645224102/e2e/$newMdName", "", "", $newMdName)
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e6c541b57a6ee7c8c57c774dde5586ca99fdcf3d/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/$zhXlf", "", "", $zhXlf)
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/bea602d80ca40d9304fa3e04a2adf05645224102/e2e/$newFile2MdName", "", "", $newFile2MdName)
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e6c541b57a6ee7c8c57c774dde5586ca99fdcf3d/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/$zhXlf", "", "", $zhXlf)
Set-HyperlinkLook $ws2.Range("A2")
Set-HyperlinkLook $ws2.Range("D2")
Set-HyperlinkLook $ws2.Range("A3")
Set-HyperlinkLook $ws2.Range("D3")

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("B2").Value2 = ".md"
$ws3.Range("C2").Value2 = "Ready for handoff"
$ws3.Range("D2").Value2 = $deXlf
$ws3.Range("E2").Value2 = $newHandoffDate
$ws3.Range("H2").Value2 = "0001-01-01 00:00:00"
$ws3.Range("J2").Value2 = "Include"
$ws3.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws3.Range("A3").Value2 = $newFile2MdName
$ws3.Range("B3").Value2 = ".md"
$ws3.Range("C3").Value2 = "Ready for handoff"
$ws3.Range("D3").Value2 = $deXlf
$ws3.Range("E3").Value2 = $newHandoffDate
$ws3.Range("H3").Value2 = "0001-01-01 00:00:00"
$ws3.Range("J3").Value2 = "Include"
$ws3.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws3.Range("A2").Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/bea602d80ca40d9304fa3e04a2adf05645224102/e2e/$newMdName", "", "", $newMdName)
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ff1aa1d73685b1597f2a7326127641e21fafbd11/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/$deXlf", "", "", $deXlf)
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/bea602d80ca40d9304fa3e04a2adf05645224102/e2e/$newFile2MdName", "", "", $newFile2MdName)
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ff1aa1d73685b1597f2a7326127641e21fafbd11/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/$deXlf", "", "", $deXlf)
Set-HyperlinkLook $ws3.Range("A2")
Set-HyperlinkLook $ws3.Range("D2")
Set-HyperlinkLook $ws3.Range("A3")
Set-HyperlinkLook $ws3.Range("D3")
